$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q0)
$ws.Range("B3").Value = 0.2704831801860494
$ws.Range("C3").Value = 1.168517949126625
$ws.Range("D3").Value = 3.573224508574976
$ws.Range("E3").Value = 1.890297465632057
$ws.Range("F3").Value = 1.877375648438149
$ws.Range("G3").Value = 144

# Row 4 (Q1)
$ws.Range("B4").Value = 0.3000337968429742
$ws.Range("C4").Value = 1.380885400760346
$ws.Range("D4").Value = 8.048481127830012
$ws.Range("E4").Value = 2.836984513145959
$ws.Range("F4").Value = 2.841443446715117
$ws.Range("G4").Value = 70
